# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 2-4, column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 1031
$wsExhibition.Range("F3").Value = 2164
$wsExhibition.Range("F4").Value = 11

# Sheet "全部类型" (All Types) - rows 4-6, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1031
$wsAll.Range("F5").Value = 2164
$wsAll.Range("F6").Value = 11
